$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (strikeouts) column values replacing the old Strike# values.
$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 3
    6  = 1
    7  = 0
    8  = 4
    9  = 1
    10 = 3
    11 = 0
    12 = 3
    13 = 1
    14 = 2
    16 = 0
    17 = 1
    18 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
